# #5: property aircraft done
# The "property_category" column on the 建物 (building) sheet was still
# tagged "land" (copy/paste leftover from the 土地 sheet); same for the
# 汽車 (car) sheet. Fix the category labels to match their own sheet.

$wb = $excel.ActiveWorkbook

# 建物 (building) sheet — column I (property_category) for every data row
# should read "building" instead of "land".
$wsBuilding = $wb.Worksheets.Item("建物")
$lastRowBuilding = $wsBuilding.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRowBuilding; $r++) {
    if ($wsBuilding.Cells.Item($r, 9).Text -eq "land") {
        $wsBuilding.Cells.Item($r, 9).Value = "building"
    }
}

# 汽車 (car) sheet — column H (property_category) should read "car"
# instead of "land".
$wsCar = $wb.Worksheets.Item("汽車")
$lastRowCar = $wsCar.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRowCar; $r++) {
    if ($wsCar.Cells.Item($r, 8).Text -eq "land") {
        $wsCar.Cells.Item($r, 8).Value = "car"
    }
}
